# Replicates the "Send email with data from excel" commit:
#   - columns A, B and C were manually widened so the full product name,
#     the full product link and the price are all visible (the data is
#     then copied/read out to build the outgoing email)
#   - the workbook window was scrolled down one row and minimized when it
#     was last saved (sheetView/topLeftCell + workbookView/minimized)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------------
# Target stored (OOXML) widths are 134.85546875 / 144.5703125 / 24.85546875
# "characters". Feed ColumnWidth the character-unit values that land
# closest to those saved widths.
$ws.Columns.Item(1).ColumnWidth = 134
$ws.Columns.Item(2).ColumnWidth = 143.666666666667
$ws.Columns.Item(3).ColumnWidth = 24

# --- Window / view state ---------------------------------------------------
# Scroll so row 2 becomes the top-left visible row, then minimize the
# window, matching the session state the file was saved in.
[void]$ws.Range("A2").Select
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$win.WindowState = -4140
$excel.WindowState = -4140
